# Updated cryptos list with latest prices/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.808.16"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "'1.764.36"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'327.20"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4475"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").Value = "'0.3536"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "'0.07425"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'42.03"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'1.090"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "'20.76"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "'6.007"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'7.183"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "'1.763.33"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'92.93"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "'0.06428"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").Value = "'5.768"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'27.849.59"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'2.116"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "'162.85"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'20.14"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "'1.971.25"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'2.155"
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("D30").Value = "'124.87"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "'1.092"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").Value = "'0.09119"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").Value = "'3.652"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'5.541"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'11.79"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").Value = "'0.02285"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'0.06094"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Value = "'0.2088"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.942"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6275"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'1.180"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").Value = "'7.889"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "'13.20"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'3.733"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "'0.5833"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "'121.84"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").Value = "'0.06899"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'1.133"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "'72.63"
$ws.Range("E51").Value = "  +0.48%  "
